$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.132.05"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "3.129.06"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'579.43"
$ws.Range("D6").Value = "'177.35"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.129.05"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").Value = "'6.41"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "'36.39"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "3.651.05"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "67.057.22"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'17.06"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "3.133.67"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'490.22"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'7.84"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").Value = "'84.03"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'12.86"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "'10.29"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").Value = "'28.18"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "0.0₃0945"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D36").Value = "'48.66"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("D38").Value = "'0.947"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'49.46"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.312"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("D45").Value = "2.813.03"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "'376.06"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'0.0349"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "'134.88"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D50").Value = "'24.92"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +1.73%  "
